# 17.1.1 Total government revenue as a proportion of GDP, by source
# - Refresh the 2019 (P) / 2020 (Q) figures with revised data.
# - Drop the 2021 (R) / 2022 (S) columns entirely (not yet available).
# - Leave the selection where the source file had it after the edit (N13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Revenues, total
$ws.Range("P5").Value = 27
$ws.Range("Q5").Value = 25.3

# Row 6 - Tax revenues
$ws.Range("P6").Value = 19.6
$ws.Range("Q6").Value = 17.8

# Row 8 - Non-tax revenues
$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 2

# Row 9 - Revenues from the sale of non-financial assets
$ws.Range("P9").Value = 5.2
$ws.Range("Q9").Value = 5.5

# Remove the now-unused 2021/2022 columns (R:S) for the whole table
$null = $ws.Range("R1:S12").Delete()

# Match the saved selection state of the edited workbook
$null = $ws.Range("N13").Select()
